$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B2 = New-Object 'object[,]' 24,4
$arr_B2[0,0] = 16.75986665991057
$arr_B2[0,1] = 6.796856553480716
$arr_B2[0,2] = 6.039812658324104
$arr_B2[0,3] = 11.36589656093263
$arr_B2[1,0] = 16.59878741259503
$arr_B2[1,1] = 6.692520845230397
$arr_B2[1,2] = 5.930605188641134
$arr_B2[1,3] = 11.38427296345069
$arr_B2[2,0] = 16.5035426900284
$arr_B2[2,1] = 6.62647670951273
$arr_B2[2,2] = 5.864366657442162
$arr_B2[2,3] = 11.39664077983394
$arr_B2[3,0] = 16.46569172169145
$arr_B2[3,1] = 6.599077209915897
$arr_B2[3,2] = 5.837618580578981
$arr_B2[3,3] = 11.40195387381332
$arr_B2[4,0] = 16.45946589134147
$arr_B2[4,1] = 6.594498561040239
$arr_B2[4,2] = 5.833193016972498
$arr_B2[4,3] = 11.40285261483527
$arr_B2[5,0] = 16.50302826934411
$arr_B2[5,1] = 6.626109141357361
$arr_B2[5,2] = 5.864004880539332
$arr_B2[5,3] = 11.39671132778077
$arr_B2[6,0] = 16.70359300771579
$arr_B2[6,1] = 6.761300106481924
$arr_B2[6,2] = 6.002012488235305
$arr_B2[6,3] = 11.37200785230183
$arr_B2[7,0] = 17.12399333958425
$arr_B2[7,1] = 7.010173732055011
$arr_B2[7,2] = 6.277380088621245
$arr_B2[7,3] = 11.33215422387745
$arr_B2[8,0] = 17.44668414877539
$arr_B2[8,1] = 7.182478707943786
$arr_B2[8,2] = 6.480297210681603
$arr_B2[8,3] = 11.30808844030197
$arr_B2[9,0] = 17.59589278638379
$arr_B2[9,1] = 7.258449930946891
$arr_B2[9,2] = 6.572277137863535
$arr_B2[9,3] = 11.29826788384732
$arr_B2[10,0] = 17.65269109312089
$arr_B2[10,1] = 7.28686164198607
$arr_B2[10,2] = 6.607023635579486
$arr_B2[10,3] = 11.29471077656167
$arr_B2[11,0] = 17.64044619880008
$arr_B2[11,1] = 7.280758703335305
$arr_B2[11,2] = 6.599544677935659
$arr_B2[11,3] = 11.29546967616882
$arr_B2[12,0] = 17.60055997257844
$arr_B2[12,1] = 7.260794589645325
$arr_B2[12,2] = 6.575137649292663
$arr_B2[12,3] = 11.29797199932002
$arr_B2[13,0] = 17.57616557102113
$arr_B2[13,1] = 7.248519215719678
$arr_B2[13,2] = 6.560175587854043
$arr_B2[13,3] = 11.29952579529203
$arr_B2[14,0] = 17.43697724716978
$arr_B2[14,1] = 7.177464404199895
$arr_B2[14,2] = 6.474276292120774
$arr_B2[14,3] = 11.30875288692921
$arr_B2[15,0] = 17.35217196566087
$arr_B2[15,1] = 7.133249565884342
$arr_B2[15,2] = 6.421468882429076
$arr_B2[15,3] = 11.31470183655848
$arr_B2[16,0] = 17.30362458730052
$arr_B2[16,1] = 7.107591983771496
$arr_B2[16,2] = 6.391066696126457
$arr_B2[16,3] = 11.31822962334282
$arr_B2[17,0] = 17.28722840050262
$arr_B2[17,1] = 7.098866207700355
$arr_B2[17,2] = 6.380769257093681
$arr_B2[17,3] = 11.31944230608288
$arr_B2[18,0] = 17.36117615911197
$arr_B2[18,1] = 7.137979815894485
$arr_B2[18,2] = 6.427093572305
$arr_B2[18,3] = 11.31405758142096
$arr_B2[19,0] = 17.61226789104132
$arr_B2[19,1] = 7.26666829965866
$arr_B2[19,2] = 6.582309160623572
$arr_B2[19,3] = 11.29723262002185
$arr_B2[20,0] = 17.77807393259218
$arr_B2[20,1] = 7.348688797477474
$arr_B2[20,2] = 6.683244148597494
$arr_B2[20,3] = 11.28717902124448
$arr_B2[21,0] = 17.68944113851592
$arr_B2[21,1] = 7.305106906083347
$arr_B2[21,2] = 6.629431486107236
$arr_B2[21,3] = 11.29245869581676
$arr_B2[22,0] = 17.35710471048342
$arr_B2[22,1] = 7.135842010098163
$arr_B2[22,2] = 6.424550779589146
$arr_B2[22,3] = 11.31434851381089
$arr_B2[23,0] = 17.00763708154662
$arr_B2[23,1] = 6.944653608597033
$arr_B2[23,2] = 6.202618001264091
$arr_B2[23,3] = 11.34201832965657
$ws.Range("B2:E25").Value = $arr_B2

$arr_G2 = New-Object 'object[,]' 24,2
$arr_G2[0,0] = 54.26803686999935
$arr_G2[0,1] = 20.54694662158878
$arr_G2[1,0] = 54.12244448164898
$arr_G2[1,1] = 20.56566637289979
$arr_G2[2,0] = 54.04496650767897
$arr_G2[2,1] = 20.58034264248835
$arr_G2[3,0] = 54.01640635441598
$arr_G2[3,1] = 20.5871230551092
$arr_G2[4,0] = 54.01184633988723
$arr_G2[4,1] = 20.58829722032272
$arr_G2[5,0] = 54.04456911788983
$arr_G2[5,1] = 20.58043084832555
$arr_G2[6,0] = 54.21537407248268
$arr_G2[6,1] = 20.55274046267525
$arr_G2[7,0] = 54.64402268618487
$arr_G2[7,1] = 20.52370689326677
$arr_G2[8,0] = 55.01474253059116
$arr_G2[8,1] = 20.51779615268609
$arr_G2[9,0] = 55.19516138535482
$arr_G2[9,1] = 20.51845539340756
$arr_G2[10,0] = 55.26514082828261
$arr_G2[10,1] = 20.51918604444164
$arr_G2[11,0] = 55.24999631194994
$arr_G2[11,1] = 20.51900730178551
$arr_G2[12,0] = 55.20088560117782
$arr_G2[12,1] = 20.51850586759348
$arr_G2[13,0] = 55.17101881259765
$arr_G2[13,1] = 20.51826134989102
$arr_G2[14,0] = 55.00318580363355
$arr_G2[14,1] = 20.51782040888875
$arr_G2[15,0] = 54.90321812398924
$arr_G2[15,1] = 20.51840718228749
$arr_G2[16,0] = 54.84683007911719
$arr_G2[16,1] = 20.51905985401502
$arr_G2[17,0] = 54.82792983586889
$arr_G2[17,1] = 20.5193349789911
$arr_G2[18,0] = 54.9137451413817
$arr_G2[18,1] = 20.51831210225512
$arr_G2[19,0] = 55.21526587986004
$arr_G2[19,1] = 20.51864010100569
$arr_G2[20,0] = 55.42197751737606
$arr_G2[20,1] = 20.52165781364649
$arr_G2[21,0] = 55.31078077795378
$arr_G2[21,1] = 20.5197908966161
$arr_G2[22,0] = 54.90898249527364
$arr_G2[22,1] = 20.51835410560165
$arr_G2[23,0] = 54.51815580474612
$arr_G2[23,1] = 20.52885372164463
$ws.Range("G2:H25").Value = $arr_G2

$arr_K2 = New-Object 'object[,]' 24,4
$arr_K2[0,0] = 13.11867552786184
$arr_K2[0,1] = 10.06044851133212
$arr_K2[0,2] = 16.15290065296539
$arr_K2[0,3] = 23.49160455172819
$arr_K2[1,0] = 13.00208010938124
$arr_K2[1,1] = 10.07064618614958
$arr_K2[1,2] = 16.14057881007925
$arr_K2[1,3] = 23.53911512472145
$arr_K2[2,0] = 12.93338772004742
$arr_K2[2,1] = 10.07829238183617
$arr_K2[2,2] = 16.13595105707765
$arr_K2[2,3] = 23.57018126463647
$arr_K2[3,0] = 12.90615285757174
$arr_K2[3,1] = 10.08175676762917
$arr_K2[3,2] = 16.13480638409189
$arr_K2[3,3] = 23.58331770221642
$arr_K2[4,0] = 12.90167712494734
$arr_K2[4,1] = 10.08235308194018
$arr_K2[4,2] = 16.13466113954007
$arr_K2[4,3] = 23.58552780190677
$arr_K2[5,0] = 12.93301731524442
$arr_K2[5,1] = 10.07833769236193
$arr_K2[5,2] = 16.13593261576304
$arr_K2[5,3] = 23.57035649651011
$arr_K2[6,0] = 13.07788998487029
$arr_K2[6,1] = 10.06367744350571
$arr_K2[6,2] = 16.14804388406285
$arr_K2[6,3] = 23.50759321683459
$arr_K2[7,0] = 13.38362544246086
$arr_K2[7,1] = 10.04590258697166
$arr_K2[7,2] = 16.19498327152645
$arr_K2[7,3] = 23.39953003639105
$arr_K2[8,0] = 13.61955118273304
$arr_K2[8,1] = 10.0395128742585
$arr_K2[8,2] = 16.24342014192866
$arr_K2[8,3] = 23.32926719530054
$arr_K2[9,0] = 13.72891430289673
$arr_K2[9,1] = 10.03804869490096
$arr_K2[9,2] = 16.26843627619346
$arr_K2[9,3] = 23.29928104533442
$arr_K2[10,0] = 13.77058460844127
$arr_K2[10,1] = 10.03770108306076
$arr_K2[10,2] = 16.27833321931923
$arr_K2[10,3] = 23.28821005572662
$arr_K2[11,0] = 13.76159931938258
$arr_K2[11,1] = 10.03776675835482
$arr_K2[11,2] = 16.27618296638699
$arr_K2[11,3] = 23.29058175742463
$arr_K2[12,0] = 13.73233760503219
$arr_K2[12,1] = 10.03801595490394
$arr_K2[12,2] = 16.26924203824645
$arr_K2[12,3] = 23.29836453503316
$arr_K2[13,0] = 13.71444632874568
$arr_K2[13,1] = 10.03819551302366
$arr_K2[13,2] = 16.26504556721856
$arr_K2[13,3] = 23.30316870721192
$arr_K2[14,0] = 13.61244197616874
$arr_K2[14,1] = 10.03963754011379
$arr_K2[14,2] = 16.24184487179395
$arr_K2[14,3] = 23.33126662580147
$arr_K2[15,0] = 13.55036211218653
$arr_K2[15,1] = 10.04089124384428
$arr_K2[15,2] = 16.22837237582113
$arr_K2[15,3] = 23.34900998002088
$arr_K2[16,0] = 13.51484956467893
$arr_K2[16,1] = 10.04174815891677
$arr_K2[16,2] = 16.22090435006442
$arr_K2[16,3] = 23.3594015582035
$arr_K2[17,0] = 13.50286007657415
$arr_K2[17,1] = 10.04206163944435
$arr_K2[17,2] = 16.21842420973276
$arr_K2[17,3] = 23.36295193041253
$arr_K2[18,0] = 13.55695078901004
$arr_K2[18,1] = 10.04074373238548
$arr_K2[18,2] = 16.22977749816987
$arr_K2[18,3] = 23.34710191315974
$arr_K2[19,0] = 13.74092579936936
$arr_K2[19,1] = 10.03793715129767
$arr_K2[19,2] = 16.27126929453498
$arr_K2[19,3] = 23.29607083485062
$arr_K2[20,0] = 13.86264436831708
$arr_K2[20,1] = 10.03730819277992
$arr_K2[20,2] = 16.3008546798081
$arr_K2[20,3] = 23.26437505864978
$arr_K2[21,0] = 13.79755748467258
$arr_K2[21,1] = 10.03753380811291
$arr_K2[21,2] = 16.28484030597283
$arr_K2[21,3] = 23.28114021096102
$arr_K2[22,0] = 13.55397149099692
$arr_K2[22,1] = 10.04080999819565
$arr_K2[22,2] = 16.22914137760864
$arr_K2[22,3] = 23.34796395612496
$arr_K2[23,0] = 13.29879059690994
$arr_K2[23,1] = 10.04953819364532
$arr_K2[23,2] = 16.17981989497649
$arr_K2[23,3] = 23.42715888840173
$ws.Range("K2:N25").Value = $arr_K2
